$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.756.29"
$ws.Range("E2").Value = "  +0.55%  "
$ws.Range("D3").Value = "2.470.85"
$ws.Range("E3").Value = "  -0.54%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").Value = "'316.40"
$ws.Range("E5").Value = "  +1.11%  "
$ws.Range("D6").Value = "'93.00"
$ws.Range("E6").Value = "  -0.22%  "
$ws.Range("D7").Value = "'0.549"
$ws.Range("E7").Value = "  +0.53%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("E9").Value = "  +3.47%  "
$ws.Range("E10").Value = "  +0.81%  "
$ws.Range("D11").Value = "'0.0853"
$ws.Range("E11").Value = "  +8.95%  "
$ws.Range("E12").Value = "  +0.09%  "
$ws.Range("D13").Value = "2.849.64"
$ws.Range("E13").Value = "  -0.67%  "
$ws.Range("D14").Value = "'6.89"
$ws.Range("E14").Value = "  +0.84%  "
$ws.Range("D15").Value = "'15.79"
$ws.Range("E15").Value = "  +0.36%  "
$ws.Range("D16").Value = "2.477.13"
$ws.Range("E16").Value = "  +0.79%  "
$ws.Range("D17").Value = "'0.783"
$ws.Range("E17").Value = "  +4.22%  "
$ws.Range("D18").Value = "41.729.00"
$ws.Range("E18").Value = "  +0.44%  "
$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").Value = "0.0₃0958"
$ws.Range("E19").Value = "  +3.23%  "
$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").Value = "'6.49"
$ws.Range("E20").Value = "  +2.44%  "
$ws.Range("D21").Value = "'11.53"
$ws.Range("E21").Value = "  +2.67%  "
$ws.Range("E22").Value = "  +0.15%  "
$ws.Range("D23").Value = "'239.93"
$ws.Range("E23").Value = "  +1.61%  "
$ws.Range("E24").Value = "  +0.76%  "
$ws.Range("E25").Value = "  +1.32%  "
$ws.Range("E26").Value = "  +0.05%  "
$ws.Range("D27").Value = "'24.72"
$ws.Range("E27").Value = "  -0.95%  "
$ws.Range("D28").Value = "'2.28"
$ws.Range("E28").Value = "  +2.19%  "
$ws.Range("E29").Value = "  +1.63%  "
$ws.Range("D30").Value = "'35.50"
$ws.Range("E30").Value = "  -2.21%  "
$ws.Range("D31").Value = "'156.03"
$ws.Range("E31").Value = "  -0.91%  "
$ws.Range("D32").Value = "'5.50"
$ws.Range("E32").Value = "  +1.28%  "
$ws.Range("E33").Value = "  +0.53%  "
$ws.Range("D34").Value = "'0.0765"
$ws.Range("E34").Value = "  +1.36%  "
$ws.Range("E35").Value = "  +2.40%  "
$ws.Range("D36").Value = "'17.48"
$ws.Range("E36").Value = "  -3.48%  "
$ws.Range("E37").Value = "  -1.85%  "
$ws.Range("E38").Value = "  +1.04%  "
$ws.Range("E39").Value = "  -2.46%  "
$ws.Range("E40").Value = "  -2.19%  "
$ws.Range("D41").Value = "'3.99"
$ws.Range("E41").Value = "  -3.90%  "
$ws.Range("E42").Value = "  -0.12%  "
$ws.Range("D43").Value = "1.970.03"
$ws.Range("E43").Value = "  +0.39%  "
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").Value = "'19.04"
$ws.Range("E44").Value = "  -3.79%  "
$ws.Range("B45").Value = "VeChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D45").Value = "'0.0285"
$ws.Range("E45").Value = "  +0.15%  "
$ws.Range("E46").Value = "  -1.09%  "
$ws.Range("D47").Value = "'9.05"
$ws.Range("E47").Value = "  +2.05%  "
$ws.Range("D48").Value = "2.702.81"
$ws.Range("D49").Value = "'97.33"
$ws.Range("E49").Value = "  +0.85%  "
$ws.Range("D50").Value = "'67.12"
$ws.Range("E50").Value = "  -0.56%  "
$ws.Range("D51").Value = "'52.79"
$ws.Range("E51").Value = "  +4.64%  "
